$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 30 and row 31 for columns
# A, D, E, F, G, H, Q, R, Z, AB (as described by the diff), while
# column B receives two brand new values (not swapped).

$row30 = @{
    A  = 112415012
    B  = 90826
    D  = "LC"
    E  = 4366
    F  = "Skarp dropptaggsvamp"
    G  = "Hydnellum peckii"
    H  = "Banker"
    Q  = 381335
    R  = 6862894
    Z  = "11:15"
    AB = "11:15"
}

$row31 = @{
    A  = 112414988
    B  = 90837
    D  = "NT"
    E  = 5966
    F  = "Motaggsvamp"
    G  = "Sarcodon squamosus"
    H  = "(Schaeff.) Quél."
    Q  = 381294
    R  = 6862860
    Z  = "11:11"
    AB = "11:11"
}

foreach ($col in $row30.Keys) {
    $ws.Range("$col" + "30").Value = $row30[$col]
}

foreach ($col in $row31.Keys) {
    $ws.Range("$col" + "31").Value = $row31[$col]
}
